$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new test case (TU08) on row 9 ------------------------------
$ws.Range("A9").Value = "TU08"
$ws.Range("B9").Value = "Check response when passing all previous test cases in an array"
$ws.Range("C9").Value = "./weather New York, 10005, Tokyo, São Paulo, Pluto, Juptior"
$ws.Range("D9").Value = "Enter all the previous arguments in array"
$ws.Range("E9").Value = "As expected"
$ws.Range("F9").Value = "Pass"

# Match the look of the other data rows: centered text, thin box border,
# regular body font, wrapped text in the middle (longer) columns.
$dataRow = $ws.Range("A9:F9")
$dataRow.HorizontalAlignment = -4108
$dataRow.VerticalAlignment = -4108
$dataRow.Borders.LineStyle = 1
$dataRow.Borders.Weight = 2
$dataRow.Font.Name = "Calibri"
$dataRow.Font.Size = 11
$dataRow.Font.Bold = $false
$dataRow.Font.Color = 0
$dataRow.WrapText = $false

$ws.Range("B9:D9").WrapText = $true

$ws.Rows(9).RowHeight = 38.25

# --- Drop the old highlight fill that used to mark row 6 ----------------
$noWrapRow6 = @("A6", "C6", "F6")
foreach ($addr in $noWrapRow6) {
    $r = $ws.Range($addr)
    $r.Interior.Pattern = -4142
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.WrapText = $false
    $r.Borders.LineStyle = 1
    $r.Borders.Weight = 2
}
$wrapRow6 = @("B6", "D6", "E6")
foreach ($addr in $wrapRow6) {
    $r = $ws.Range($addr)
    $r.Interior.Pattern = -4142
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.WrapText = $true
    $r.Borders.LineStyle = 1
    $r.Borders.Weight = 2
}

# --- Leave the selection where data entry finished -----------------------
$ws.Range("E9").Select()
